# Daily attendance processing - 2026-01-18 23:01:22
# For every cell in column G ("Recorded By") whose value is the combined
# list "dnasr281@gmail.com, System", reorder it to "System, dnasr281@gmail.com".
# Leave every other cell (including other "Recorded By" values such as a lone
# "dnasr281@gmail.com") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count
$lastRowIndex = $ws.UsedRange.Row + $lastRow - 1

for ($r = 1; $r -le $lastRowIndex; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
